# Apply updated "dSF" (column F) values for the listed rows, reflecting the
# repulled data / mean calculation described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -4
    3  = -3
    13 = -3
    14 = 2
    16 = -3
    17 = 1
    25 = 0
    27 = 1
    31 = 1
    41 = 4
    44 = 0
    46 = 0
    49 = -1
    51 = 4
    52 = 0
    56 = 1
    57 = 1
    60 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
